# Normalize the "Recorded By" (column G) entries so that any "System" /
# "system" token in the comma-separated list is moved to the front of the
# list. When two such tokens exist in one cell (one written as lowercase
# "system" and one as "System"), their text (and therefore their casing)
# is swapped between the first and last matching slot.
#
# Example:
#   "dnasr281@gmail.com, System"              -> "System, dnasr281@gmail.com"
#   "system, backup@backdoor.com, System"     -> "System, backup@backdoor.com, system"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rowCount = $used.Rows.Count

for ($r = 2; $r -le $rowCount; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $orig = $cell.Value2

    if ($null -eq $orig) { continue }
    if ($orig -eq "") { continue }

    # Split on commas and trim whitespace from each token.
    $rawParts = $orig.Split(",")
    $parts = @()
    foreach ($p in $rawParts) {
        $parts += $p.Trim()
    }

    # Locate every token that equals "system" (case-insensitively).
    $idxs = @()
    for ($i = 0; $i -lt $parts.Count; $i++) {
        if ($parts[$i].ToLower() -eq "system") {
            $idxs += $i
        }
    }

    if ($idxs.Count -eq 0) {
        # No "System" token present - leave the cell untouched.
        continue
    } elseif ($idxs.Count -eq 1) {
        # Single match: move it to the front of the list.
        $i = $idxs[0]
        $val = $parts[$i]
        $newParts = @($val)
        for ($k = 0; $k -lt $parts.Count; $k++) {
            if ($k -ne $i) { $newParts += $parts[$k] }
        }
        $parts = $newParts
    } else {
        # Multiple matches: swap their text between symmetric slots
        # (first <-> last, etc.), which also moves the first one to
        # position 0.
        $vals = @()
        foreach ($i in $idxs) { $vals += $parts[$i] }
        $n = $vals.Count
        for ($k = 0; $k -lt $idxs.Count; $k++) {
            $parts[$idxs[$k]] = $vals[$n - 1 - $k]
        }
    }

    $newVal = $parts -join ", "
    $cell.Value2 = $newVal
}
